$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update Objetivos (row 10) B/C with the objectives essay text ---
$ws.Cells.Item(10,2).Value = "Levar os alunos a`n1. conhecer mecanismos discursivos para a elaboração de monografias, trabalhos de fim de curso  e projetos de pesquisa e`n2. dominar procedimentos para apresentação de trabalho em congressos e elaboração de artigos para publicação."
$ws.Cells.Item(10,3).Value = "Levar os alunos a`n1. conhecer mecanismos discursivos para a elaboração de monografias, trabalhos de fim de curso  e projetos de pesquisa e`n2. dominar procedimentos para apresentação de trabalho em congressos e elaboração de artigos para publicação."

# --- 2. Insert a new blank row at 13, pushing old rows 13-21 down to 14-22 ---
$ws.Rows.Item(13).Insert()

# Remove the stray A13 cell created by the insert (final layout has no A13 cell)
$ws.Range("A13").Clear()

# Copy B/C formatting down from the row below (old row 13, now row 14) into new row 13
$ws.Range("B14").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C13").PasteSpecial(-4122)

# New row 13: professor name (moved from old row 10)
$ws.Cells.Item(13,2).Value = "5840514 - Graziela Zamponi"
$ws.Cells.Item(13,3).Value = "5840514 - Graziela Zamponi"

# --- 3. Row 14 ("Programa resumido:") B/C: "Semestral" -> short syllabus text ---
$ws.Cells.Item(14,2).Value = "Elaboração de textos acadêmicos: projetos de pesquisa, monografias, TCC para publicação e apresentação em eventos científicos."
$ws.Cells.Item(14,3).Value = "Elaboração de textos acadêmicos: projetos de pesquisa, monografias, TCC para publicação e apresentação em eventos científicos."

# --- 4. Row 16 ("Programa:") B/C: "01/01/2012" -> full syllabus text ---
$ws.Cells.Item(16,2).Value = "1.  Projeto de pesquisa, monografias e TCC`n1.1 - Elementos constitutivos`n1.2 - Etapas de elaboração`n2.  O processo de divulgação do trabalho científico`n2.1 - O percurso de divulgação da pesquisa`n2.2 - Procedimentos para a apresentação pública de trabalhos científicos `n2.3 - Etapas de produção de artigo para publicação"
$ws.Cells.Item(16,3).Value = "1.  Projeto de pesquisa, monografias e TCC`n1.1 - Elementos constitutivos`n1.2 - Etapas de elaboração`n2.  O processo de divulgação do trabalho científico`n2.1 - O percurso de divulgação da pesquisa`n2.2 - Procedimentos para a apresentação pública de trabalhos científicos `n2.3 - Etapas de produção de artigo para publicação"

# --- 5. Row 19 ("Método:") B/C: professor name -> method text ---
$ws.Cells.Item(19,2).Value = "Trabalhos em grupos, Apresentação Oral e Redação de Projetos  Científicos"
$ws.Cells.Item(19,3).Value = "Trabalhos em grupos, Apresentação Oral e Redação de Projetos  Científicos"

# --- 6. Row 20 ("Critério:") B/C: method text -> criteria text ---
$ws.Cells.Item(20,2).Value = "MF = (N1+N2)/2, onde N1  =  Apresentação oral de um projeto de pesquisa e/ou pesquisa já finalizada e N 2 = Elaboração escrita de um projeto de pesquisa e seu desenvolvimento em artigo e/ou monografia/TCC"
$ws.Cells.Item(20,3).Value = "MF = (N1+N2)/2, onde N1  =  Apresentação oral de um projeto de pesquisa e/ou pesquisa já finalizada e N 2 = Elaboração escrita de um projeto de pesquisa e seu desenvolvimento em artigo e/ou monografia/TCC"

# --- 7. Row 21 ("Norma de recuperação:") B/C: criteria text -> recovery text ---
$ws.Cells.Item(21,2).Value = "Ao aluno que não alcançar a média 5,0 (cinco) no final do período letivo será dada uma nova oportunidade para a reelaboração dos trabalhos."
$ws.Cells.Item(21,3).Value = "Ao aluno que não alcançar a média 5,0 (cinco) no final do período letivo será dada uma nova oportunidade para a reelaboração dos trabalhos."

# --- 8. Row 22 ("Bibliografia:") B/C: recovery text -> bibliography text ---
$ws.Cells.Item(22,2).Value = "ANDRADE, Maria Margarida de. Introdução à metodologia do trabalho científico. 10.ed. São Paulo: Atlas, 2010.`nFARACO, Carlos Alberto; TEZZA, Cristóvão. Oficina de texto.  6. ed. Petrópolis: Vozes, 2008.`nMACHADO, A.R (coord.); LOUSADA, E.; ABREU-TARDELLI, L. S. Resenha. São Paulo: Parábola Editorial, 2004 a.`n______.  Resumo.  São Paulo: Parábola Editorial, 2004 b.`n______.  Planejar gêneros acadêmicos.  São Paulo: Parábola Editorial, 2009.`nSERAFINI, Maria José. Como escrever textos. 5.ed. São Paulo: Globo, 1992.`nSEVERINO, Antonio Joaquim. Metodologia do trabalho científico. 23.ed. São Paulo: Cortez, 2009."
$ws.Cells.Item(22,3).Value = "ANDRADE, Maria Margarida de. Introdução à metodologia do trabalho científico. 10.ed. São Paulo: Atlas, 2010.`nFARACO, Carlos Alberto; TEZZA, Cristóvão. Oficina de texto.  6. ed. Petrópolis: Vozes, 2008.`nMACHADO, A.R (coord.); LOUSADA, E.; ABREU-TARDELLI, L. S. Resenha. São Paulo: Parábola Editorial, 2004 a.`n______.  Resumo.  São Paulo: Parábola Editorial, 2004 b.`n______.  Planejar gêneros acadêmicos.  São Paulo: Parábola Editorial, 2009.`nSERAFINI, Maria José. Como escrever textos. 5.ed. São Paulo: Globo, 1992.`nSEVERINO, Antonio Joaquim. Metodologia do trabalho científico. 23.ed. São Paulo: Cortez, 2009."

